# Commit: "commit just before python application webinar on 2024-0527"
#
# The author typed the workbook's own file name into cell A1 of Sheet1
# (presumably as an on-screen label while preparing a demo for the
# webinar), then left the selection sitting on A2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "df.xlsx"

# Leave the cursor on A2, matching the saved selection in the workbook.
$ws.Range("A2").Select()
